$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.57076378422843
$ws.Range("C2").Value = 14.1194058293842
$ws.Range("D2").Value = 0.111248575415226
$ws.Range("E2").Value = 0.913401583030347
$ws.Range("F2").Value = 0.246030682774803
$ws.Range("G2").Value = 19.1945420748638
$ws.Range("H2").Value = 0.0128177417213299
$ws.Range("I2").Value = 0.989997682315297
$ws.Range("J2").Value = 12.3247707575433
$ws.Range("K2").Value = 0.118929835753016
$ws.Range("L2").Value = 103.630604377007
$ws.Range("M2").Value = 0.0000000000000112304318168007
$ws.Range("N2").Value = 12.7324704197893
$ws.Range("O2").Value = 0.139409085877251
$ws.Range("P2").Value = 91.3317115571658
$ws.Range("Q2").Value = 0.000000000000000331109300294755
$ws.Range("R2").Value = 13.1311553081888
$ws.Range("S2").Value = 0.0433677079495495
$ws.Range("T2").Value = 302.786472447761
$ws.Range("U2").Value = 0.000000000000000000000582439346228936
$ws.Range("V2").Value = 13.2384482848279
$ws.Range("W2").Value = 0.0512028602799419
$ws.Range("X2").Value = 258.548999263892
$ws.Range("Y2").Value = 0.000000000000000000000409472792389334
$ws.Range("B3").Value = 0.00663909830126552
$ws.Range("C3").Value = 0.0108645694427629
$ws.Range("D3").Value = 0.611077902004479
$ws.Range("E3").Value = 0.555175458349891
$ws.Range("F3").Value = 0.000324290628228811
$ws.Range("G3").Value = 0.0110134704575029
$ws.Range("H3").Value = 0.0294449083492922
$ws.Range("I3").Value = 0.977105801289737
$ws.Range("J3").Value = -0.002984426843591
$ws.Range("K3").Value = 0.0146542424886238
$ws.Range("L3").Value = -0.203656166185856
$ws.Range("M3").Value = 0.842821677206534
$ws.Range("N3").Value = -0.0211196258494885
$ws.Range("O3").Value = 0.0174140060443276
$ws.Range("P3").Value = -1.21279536688618
$ws.Range("Q3").Value = 0.253871985609885
$ws.Range("R3").Value = -0.0396143653236384
$ws.Range("S3").Value = 0.036157801559622
$ws.Range("T3").Value = -1.09559662410107
$ws.Range("U3").Value = 0.29967818948755
$ws.Range("V3").Value = 0.00236571744851758
$ws.Range("W3").Value = 0.0274407489185927
$ws.Range("X3").Value = 0.0862118397546604
$ws.Range("Y3").Value = 0.93286819975659
$ws.Range("B4").Value = 0.0000000701749077181736
$ws.Range("C4").Value = 0.0000551800848703746
$ws.Range("D4").Value = 0.00127174338138522
$ws.Range("E4").Value = 0.999036944939867
$ws.Range("F4").Value = 0.00000744174820629739
$ws.Range("G4").Value = 0.0000458559422974967
$ws.Range("H4").Value = 0.162285362233274
$ws.Range("I4").Value = 0.877723613044885
$ws.Range("J4").Value = 0.0000391726847633002
$ws.Range("K4").Value = 0.0000889862546544473
$ws.Range("L4").Value = 0.440210512459661
$ws.Range("M4").Value = 0.678989788862418
$ws.Range("N4").Value = 0.00026466731346922
$ws.Range("O4").Value = 0.000122362775130683
$ws.Range("P4").Value = 2.16297246598532
$ws.Range("Q4").Value = 0.0847195778544399
$ws.Range("R4").Value = 0.00103007263782878
$ws.Range("S4").Value = 0.000327251347516458
$ws.Range("T4").Value = 3.14764979776586
$ws.Range("U4").Value = 0.027903246276479
$ws.Range("B5").Value = 0.0421247656651114
$ws.Range("C5").Value = 0.0141813212878226
$ws.Range("D5").Value = 2.97044011697864
$ws.Range("E5").Value = 0.0128290743308983
$ws.Range("F5").Value = 0.0474776350954371
$ws.Range("G5").Value = 0.0128324387643661
$ws.Range("H5").Value = 3.6998138831783
$ws.Range("I5").Value = 0.00354663838588354
$ws.Range("J5").Value = 0.0416740030055535
$ws.Range("K5").Value = 0.0141751734890369
$ws.Range("L5").Value = 2.93992895662155
$ws.Range("M5").Value = 0.0135373326114084
$ws.Range("N5").Value = 0.0549480034140824
$ws.Range("O5").Value = 0.0148418197145689
$ws.Range("P5").Value = 3.702241670551
$ws.Range("Q5").Value = 0.00352127575995369
$ws.Range("R5").Value = 0.0726931130391353
$ws.Range("S5").Value = 0.0325251218923231
$ws.Range("T5").Value = 2.23498356992455
$ws.Range("U5").Value = 0.04725908186327
$ws.Range("B6").Value = 0.0554279838706064
$ws.Range("C6").Value = 0.020741759376999
$ws.Range("D6").Value = 2.67228940723668
$ws.Range("E6").Value = 0.0221528553925902
$ws.Range("F6").Value = 0.0650004224121209
$ws.Range("G6").Value = 0.0152060510957056
$ws.Range("H6").Value = 4.27464185165586
$ws.Range("I6").Value = 0.00139825630931756
$ws.Range("J6").Value = 0.052353431481137
$ws.Range("K6").Value = 0.0166295426366686
$ws.Range("L6").Value = 3.14821836204299
$ws.Range("M6").Value = 0.00956952556547372
$ws.Range("N6").Value = 0.0665903136180718
$ws.Range("O6").Value = 0.0165040461635852
$ws.Range("P6").Value = 4.03478716419236
$ws.Range("Q6").Value = 0.00207603447551493
$ws.Range("R6").Value = 0.0755273086571043
$ws.Range("S6").Value = 0.0333093924621801
$ws.Range("T6").Value = 2.26744779998191
$ws.Range("U6").Value = 0.0451943850376941
$ws.Range("B7").Value = 0.0519565023917307
$ws.Range("C7").Value = 0.0203842160861286
$ws.Range("D7").Value = 2.54885947893218
$ws.Range("E7").Value = 0.0269520763783054
$ws.Range("F7").Value = 0.0531647647393375
$ws.Range("G7").Value = 0.016260037180973
$ws.Range("H7").Value = 3.26965825155364
$ws.Range("I7").Value = 0.00740098886642791
$ws.Range("J7").Value = 0.0459627671644648
$ws.Range("K7").Value = 0.0167922816136395
$ws.Range("L7").Value = 2.73713651438121
$ws.Range("M7").Value = 0.0192111450250149
$ws.Range("N7").Value = 0.0576903834087564
$ws.Range("O7").Value = 0.0214036678152967
$ws.Range("P7").Value = 2.69535034399696
$ws.Range("Q7").Value = 0.0206784860087683
$ws.Range("R7").Value = 0.0683581566843055
$ws.Range("S7").Value = 0.0336684100720061
$ws.Range("T7").Value = 2.03033515803416
$ws.Range("U7").Value = 0.0670321520122528
$ws.Range("B8").Value = 0.0199160971591002
$ws.Range("C8").Value = 0.0320440543036722
$ws.Range("D8").Value = 0.621522388220951
$ws.Range("E8").Value = 0.549144066610772
$ws.Range("F8").Value = 0.0290611727622559
$ws.Range("G8").Value = 0.0273904891110881
$ws.Range("H8").Value = 1.06099502803298
$ws.Range("I8").Value = 0.315460941475024
$ws.Range("J8").Value = 0.0237198967432038
$ws.Range("K8").Value = 0.0296197781120893
$ws.Range("L8").Value = 0.800812776295664
$ws.Range("M8").Value = 0.443204228069854
$ws.Range("N8").Value = 0.0427484632766967
$ws.Range("O8").Value = 0.0389152188231314
$ws.Range("P8").Value = 1.09850245146988
$ws.Range("Q8").Value = 0.299572569318787
$ws.Range("R8").Value = 0.0688335137431654
$ws.Range("S8").Value = 0.057446636717022
$ws.Range("T8").Value = 1.1982166002552
$ws.Range("U8").Value = 0.260401500349019
$ws.Range("B9").Value = 0.047515270416094
$ws.Range("C9").Value = 0.0388792226448756
$ws.Range("D9").Value = 1.22212501134862
$ws.Range("E9").Value = 0.250067235151742
$ws.Range("F9").Value = 0.0600444850426103
$ws.Range("G9").Value = 0.0361355707430333
$ws.Range("H9").Value = 1.66164485043277
$ws.Range("I9").Value = 0.128066738391064
$ws.Range("J9").Value = 0.0520380131872119
$ws.Range("K9").Value = 0.0376941649917983
$ws.Range("L9").Value = 1.3805323237306
$ws.Range("M9").Value = 0.19797350401341
$ws.Range("N9").Value = 0.0884218786583162
$ws.Range("O9").Value = 0.0444400486803553
$ws.Range("P9").Value = 1.98968905939572
$ws.Range("Q9").Value = 0.0751260845321074
$ws.Range("R9").Value = 0.110847239159286
$ws.Range("S9").Value = 0.0666218119360313
$ws.Range("T9").Value = 1.66382804577154
$ws.Range("U9").Value = 0.12757953622843
$ws.Range("B10").Value = 0.0899242567470835
$ws.Range("C10").Value = 0.0387128766263198
$ws.Range("D10").Value = 2.32285132451114
$ws.Range("E10").Value = 0.0431478335647195
$ws.Range("F10").Value = 0.0996232485015071
$ws.Range("G10").Value = 0.0363194501505454
$ws.Range("H10").Value = 2.74297237674483
$ws.Range("I10").Value = 0.0211928246777691
$ws.Range("J10").Value = 0.0874221536964988
$ws.Range("K10").Value = 0.0353203038700513
$ws.Range("L10").Value = 2.47512462005248
$ws.Range("M10").Value = 0.0333744345969894
$ws.Range("N10").Value = 0.128060637192092
$ws.Range("O10").Value = 0.0455201678775246
$ws.Range("P10").Value = 2.81327251552868
$ws.Range("Q10").Value = 0.018774237958734
$ws.Range("R10").Value = 0.127720992340788
$ws.Range("S10").Value = 0.0756719502733686
$ws.Range("T10").Value = 1.6878247736366
$ws.Range("U10").Value = 0.122994869045303
$ws.Range("B11").Value = 0.0636338824414552
$ws.Range("C11").Value = 0.102918140367737
$ws.Range("D11").Value = 0.618296076999497
$ws.Range("E11").Value = 0.553737448103284
$ws.Range("F11").Value = 0.0623356416624499
$ws.Range("G11").Value = 0.116390164551355
$ws.Range("H11").Value = 0.535574822002638
$ws.Range("I11").Value = 0.606960236644958
$ws.Range("J11").Value = 0.0485895314434686
$ws.Range("K11").Value = 0.113668840695944
$ws.Range("L11").Value = 0.427465707805026
$ws.Range("M11").Value = 0.680427662871396
$ws.Range("N11").Value = 0.100784355706279
$ws.Range("O11").Value = 0.0954187088377144
$ws.Range("P11").Value = 1.05623265011572
$ws.Range("Q11").Value = 0.321998441943263
$ws.Range("R11").Value = 0.135617470574639
$ws.Range("S11").Value = 0.134103505778444
$ws.Range("T11").Value = 1.01128952436707
$ws.Range("U11").Value = 0.34176162932562
$ws.Range("B12").Value = 0.0582993286044237
$ws.Range("C12").Value = 0.0305037931463405
$ws.Range("D12").Value = 1.91121570765758
$ws.Range("E12").Value = 0.0794103219085271
$ws.Range("F12").Value = 0.0489653003937446
$ws.Range("G12").Value = 0.0314174236728701
$ws.Range("H12").Value = 1.55853964677656
$ws.Range("I12").Value = 0.144180805788516
$ws.Range("J12").Value = 0.0149017544776975
$ws.Range("K12").Value = 0.0268041057679214
$ws.Range("L12").Value = 0.555950443067256
$ws.Range("M12").Value = 0.587979209956025
$ws.Range("N12").Value = -0.0360561891948333
$ws.Range("O12").Value = 0.022474299055471
$ws.Range("P12").Value = -1.60432986612127
$ws.Range("Q12").Value = 0.137680814631377
$ws.Range("B13").Value = 0.0611547722464341
$ws.Range("C13").Value = 0.0307513438359481
$ws.Range("D13").Value = 1.98868617165747
$ws.Range("E13").Value = 0.067700920534187
$ws.Range("F13").Value = 0.064484450772735
$ws.Range("G13").Value = 0.0375169976254326
$ws.Range("H13").Value = 1.71880627060149
$ws.Range("I13").Value = 0.108734364866926
$ws.Range("J13").Value = 0.0844250179768771
$ws.Range("K13").Value = 0.0439770614471129
$ws.Range("L13").Value = 1.91975123391105
$ws.Range("M13").Value = 0.0763859281153418
$ws.Range("N13").Value = 0.108174502188427
$ws.Range("O13").Value = 0.0439388297902726
$ws.Range("P13").Value = 2.46193407300017
$ws.Range("Q13").Value = 0.0280632444491428
$ws.Range("B14").Value = 0.0244156668464332
$ws.Range("C14").Value = 0.0247779116270791
$ws.Range("D14").Value = 0.985380334464908
$ws.Range("E14").Value = 0.343897603131135
$ws.Range("F14").Value = 0.0180133933646897
$ws.Range("G14").Value = 0.0267393062680447
$ws.Range("H14").Value = 0.673667191815552
$ws.Range("I14").Value = 0.51339013201172
$ws.Range("J14").Value = 0.00503176499693464
$ws.Range("K14").Value = 0.0339483359291381
$ws.Range("L14").Value = 0.148218310536271
$ws.Range("M14").Value = 0.88466876771673
$ws.Range("N14").Value = -0.00355268578238182
$ws.Range("O14").Value = 0.0340804108084061
$ws.Range("P14").Value = -0.104244218250607
$ws.Range("Q14").Value = 0.918699399155162
$ws.Range("B15").Value = -0.033829336129902
$ws.Range("C15").Value = 0.023285341721079
$ws.Range("D15").Value = -1.45281682077606
$ws.Range("E15").Value = 0.177652466171965
$ws.Range("F15").Value = -0.0339515318027743
$ws.Range("G15").Value = 0.024089896235186
$ws.Range("H15").Value = -1.40936812144439
$ws.Range("I15").Value = 0.189729601995759
$ws.Range("J15").Value = -0.0367995018188028
$ws.Range("K15").Value = 0.0280809059717646
$ws.Range("L15").Value = -1.31048128774067
$ws.Range("M15").Value = 0.220101670491879
$ws.Range("N15").Value = -0.107723802792713
$ws.Range("O15").Value = 0.0568034923295771
$ws.Range("P15").Value = -1.89642922247973
$ws.Range("Q15").Value = 0.0836098464245918
$ws.Range("B16").Value = -0.00802107757286733
$ws.Range("C16").Value = 0.0102019160028067
$ws.Range("D16").Value = -0.786232465613382
$ws.Range("E16").Value = 0.450800367178755
$ws.Range("F16").Value = -0.00646576057092813
$ws.Range("G16").Value = 0.0117321894523018
$ws.Range("H16").Value = -0.551112867484388
$ws.Range("I16").Value = 0.594145639178856
$ws.Range("J16").Value = 0.00295625277439457
$ws.Range("K16").Value = 0.00737599131337616
$ws.Range("L16").Value = 0.400793960946441
$ws.Range("M16").Value = 0.697278405479685
$ws.Range("N16").Value = 0.0821451383531718
$ws.Range("O16").Value = 0.0189841981604653
$ws.Range("P16").Value = 4.32702701788268
$ws.Range("Q16").Value = 0.00121939542724005
$ws.Range("B17").Value = 0.000521658317935731
$ws.Range("C17").Value = 0.0000657115166975988
$ws.Range("D17").Value = 7.93861326221365
$ws.Range("E17").Value = 0.0000349600518030601
$ws.Range("F17").Value = 0.000567831141208356
$ws.Range("G17").Value = 0.000075224759177545
$ws.Range("H17").Value = 7.54846073841412
$ws.Range("I17").Value = 0.000056191710316165
$ws.Range("J17").Value = 0.000655145611916871
$ws.Range("K17").Value = 0.000065950091907328
$ws.Range("L17").Value = 9.93396055971342
$ws.Range("M17").Value = 0.00000480169921190597
$ws.Range("N17").Value = 0.000718077639688489
$ws.Range("O17").Value = 0.0000603866664307338
$ws.Range("P17").Value = 11.8913277074527
$ws.Range("Q17").Value = 0.0000016093601132628
$ws.Range("B18").Value = -0.0223832239337787
$ws.Range("C18").Value = 0.00673226713545956
$ws.Range("D18").Value = -3.32476764266883
$ws.Range("E18").Value = 0.0111944201767327
$ws.Range("F18").Value = -0.0241533244108377
$ws.Range("G18").Value = 0.00673139199849138
$ws.Range("H18").Value = -3.5881619160273
$ws.Range("I18").Value = 0.00773344432583068
$ws.Range("J18").Value = -0.0214319661656497
$ws.Range("K18").Value = 0.00818398321943224
$ws.Range("L18").Value = -2.61876956379396
$ws.Range("M18").Value = 0.031248047823141
$ws.Range("N18").Value = -0.0244911343132044
$ws.Range("O18").Value = 0.00787064045268441
$ws.Range("P18").Value = -3.11170793030589
$ws.Range("Q18").Value = 0.0152431409973079
$ws.Range("B19").Value = 0.222949567659541
$ws.Range("C19").Value = 0.0911124447437027
$ws.Range("D19").Value = 2.44697163254365
$ws.Range("E19").Value = 0.0427453870472583
$ws.Range("F19").Value = 0.257022441762782
$ws.Range("G19").Value = 0.0981098866310113
$ws.Range("H19").Value = 2.6197404827246
$ws.Range("I19").Value = 0.0330662916632569
$ws.Range("J19").Value = 0.255557658408985
$ws.Range("K19").Value = 0.115202970155184
$ws.Range("L19").Value = 2.21832525728057
$ws.Range("M19").Value = 0.0597924314013831
$ws.Range("N19").Value = 0.275493066157687
$ws.Range("O19").Value = 0.121593986128398
$ws.Range("P19").Value = 2.26568003015197
$ws.Range("Q19").Value = 0.0569056172315191
$ws.Range("B20").Value = -0.00659639869134638
$ws.Range("C20").Value = 0.1327255899377
$ws.Range("D20").Value = -0.0496995243678531
$ws.Range("E20").Value = 0.961996334678582
$ws.Range("F20").Value = -0.012291609083111
$ws.Range("G20").Value = 0.148430666915383
$ws.Range("H20").Value = -0.0828104416597294
$ws.Range("I20").Value = 0.936731258887103
$ws.Range("J20").Value = 0.0394806044235857
$ws.Range("K20").Value = 0.131677518795238
$ws.Range("L20").Value = 0.299827979633936
$ws.Range("M20").Value = 0.774561756860222
$ws.Range("N20").Value = 0.0106668635307851
$ws.Range("O20").Value = 0.173721429288596
$ws.Range("P20").Value = 0.0614021170241736
$ws.Range("Q20").Value = 0.953060989047301
$ws.Range("B21").Value = 0.0997983889288878
$ws.Range("C21").Value = 0.0130551514691391
$ws.Range("D21").Value = 7.64436852110064
$ws.Range("E21").Value = 0.0000116464787357637
$ws.Range("F21").Value = 0.105041992013948
$ws.Range("G21").Value = 0.0171109876793328
$ws.Range("H21").Value = 6.1388619980611
$ws.Range("I21").Value = 0.0000822248607869245
$ws.Range("J21").Value = 0.134184384948773
$ws.Range("K21").Value = 0.0179285133651969
$ws.Range("L21").Value = 7.48441224408792
$ws.Range("M21").Value = 0.0000137433646471336
$ws.Range("N21").Value = 0.157391272847764
$ws.Range("O21").Value = 0.0292494251399044
$ws.Range("P21").Value = 5.38100397169988
$ws.Range("Q21").Value = 0.000250310041870651
$ws.Range("B22").Value = 0.0569629889544122
$ws.Range("C22").Value = 0.0061539965794301
$ws.Range("D22").Value = 9.25625944362929
$ws.Range("E22").Value = 0.0000014001164391205
$ws.Range("F22").Value = 0.0590268006694955
$ws.Range("G22").Value = 0.00769772351830218
$ws.Range("H22").Value = 7.66808531498343
$ws.Range("I22").Value = 0.00000880129574730221
$ws.Range("J22").Value = 0.0722778024969536
$ws.Range("K22").Value = 0.0125193690401522
$ws.Range("L22").Value = 5.7732783709102
$ws.Range("M22").Value = 0.000115690661512462
$ws.Range("N22").Value = 0.07760744693659
$ws.Range("O22").Value = 0.0142653202947853
$ws.Range("P22").Value = 5.4402877280617
$ws.Range("Q22").Value = 0.000197420307530553
$ws.Range("B23").Value = -0.0926070120758595
$ws.Range("C23").Value = 0.0112404323939888
$ws.Range("D23").Value = -8.23874107595579
$ws.Range("E23").Value = 0.00000398228520027614
$ws.Range("F23").Value = -0.0950847164114034
$ws.Range("G23").Value = 0.0116554920194514
$ws.Range("H23").Value = -8.15793243671909
$ws.Range("I23").Value = 0.00000439675832219963
$ws.Range("J23").Value = -0.105475696620902
$ws.Range("K23").Value = 0.0165771483585578
$ws.Range("L23").Value = -6.36271657461827
$ws.Range("M23").Value = 0.0000462042576293582
$ws.Range("N23").Value = -0.091251875713418
$ws.Range("O23").Value = 0.0158070007800941
$ws.Range("P23").Value = -5.77287728285129
$ws.Range("Q23").Value = 0.000110088709089378
$ws.Range("B24").Value = -0.0240102197893225
$ws.Range("C24").Value = 0.0209041844530363
$ws.Range("D24").Value = -1.14858438238833
$ws.Range("E24").Value = 0.276350605465039
$ws.Range("F24").Value = -0.0234299084832179
$ws.Range("G24").Value = 0.0220484449664152
$ws.Range("H24").Value = -1.06265582533857
$ws.Range("I24").Value = 0.311936849567052
$ws.Range("J24").Value = -0.00185509715120677
$ws.Range("K24").Value = 0.0259141445560999
$ws.Range("L24").Value = -0.071586277802487
$ws.Range("M24").Value = 0.944281805329166
$ws.Range("B25").Value = 0.0177810224630068
$ws.Range("C25").Value = 0.0506914220545293
$ws.Range("D25").Value = 0.350769849065973
$ws.Range("E25").Value = 0.737306557293768
$ws.Range("F25").Value = 0.029414055232081
$ws.Range("G25").Value = 0.0633367567907874
$ws.Range("H25").Value = 0.464407347683445
$ws.Range("I25").Value = 0.658171543322011
$ws.Range("J25").Value = 0.0428928418011134
$ws.Range("K25").Value = 0.0908387842113929
$ws.Range("L25").Value = 0.472186436371678
$ws.Range("M25").Value = 0.65300233711491
$ws.Range("B26").Value = -0.0254669466126221
$ws.Range("C26").Value = 0.00927944584244082
$ws.Range("D26").Value = -2.7444469254991
$ws.Range("E26").Value = 0.0187902819372597
$ws.Range("F26").Value = -0.0301008683076392
$ws.Range("G26").Value = 0.0119768211489444
$ws.Range("H26").Value = -2.51326023268638
$ws.Range("I26").Value = 0.0285022738337785
$ws.Range("J26").Value = -0.0133983289051608
$ws.Range("K26").Value = 0.0128387416233555
$ws.Range("L26").Value = -1.0435858356077
$ws.Range("M26").Value = 0.318793624976163
$ws.Range("B27").Value = 0.0620553266277024
$ws.Range("C27").Value = 0.0222169973124795
$ws.Range("D27").Value = 2.79314642545531
$ws.Range("E27").Value = 0.024901024566749
$ws.Range("F27").Value = 0.0838159771011283
$ws.Range("G27").Value = 0.0225915310603414
$ws.Range("H27").Value = 3.71006183145613
$ws.Range("I27").Value = 0.00653667375806786
$ws.Range("J27").Value = 0.0818956467833765
$ws.Range("K27").Value = 0.0271021290084299
$ws.Range("L27").Value = 3.02174219441961
$ws.Range("M27").Value = 0.0172794773426854
$ws.Range("B28").Value = 0.0171759396464268
$ws.Range("C28").Value = 0.0265608067388243
$ws.Range("D28").Value = 0.646664832710841
$ws.Range("E28").Value = 0.53033289992282
$ws.Range("F28").Value = 0.0235237396383613
$ws.Range("G28").Value = 0.0246147430624193
$ws.Range("H28").Value = 0.955676830698926
$ws.Range("I28").Value = 0.358604330894163
$ws.Range("J28").Value = -0.00798737548632356
$ws.Range("K28").Value = 0.0393504967198022
$ws.Range("L28").Value = -0.202980296365716
$ws.Range("M28").Value = 0.842566757686755
$ws.Range("B29").Value = 0.0107379841416726
$ws.Range("C29").Value = 0.118873990875185
$ws.Range("D29").Value = 0.0903308121702351
$ws.Range("E29").Value = 0.938263871849925
$ws.Range("F29").Value = 0.0366716056016444
$ws.Range("G29").Value = 0.129376103368041
$ws.Range("H29").Value = 0.283449606588655
$ws.Range("I29").Value = 0.810075545752507
$ws.Range("J29").Value = 0.0454580392364991
$ws.Range("K29").Value = 0.125208449884483
$ws.Range("L29").Value = 0.363058877243818
$ws.Range("M29").Value = 0.760087370365826
$ws.Range("B30").Value = -0.00292257674260236
$ws.Range("C30").Value = 0.0123871898960521
$ws.Range("D30").Value = -0.235935411269816
$ws.Range("E30").Value = 0.8172814229356
$ws.Range("F30").Value = -0.010112710445206
$ws.Range("G30").Value = 0.0134321863972948
$ws.Range("H30").Value = -0.75287150923119
$ws.Range("I30").Value = 0.465365058998382
$ws.Range("J30").Value = -0.0362419189523006
$ws.Range("K30").Value = 0.0267761969137503
$ws.Range("L30").Value = -1.3535125645005
$ws.Range("M30").Value = 0.199584238636603
$ws.Range("B31").Value = -0.0371825979264496
$ws.Range("C31").Value = 0.0757798832953121
$ws.Range("D31").Value = -0.49066581142056
$ws.Range("E31").Value = 0.645833760822474
$ws.Range("F31").Value = -0.0356370022599731
$ws.Range("G31").Value = 0.075403126507532
$ws.Range("H31").Value = -0.47261968980044
$ws.Range("I31").Value = 0.657819456550932
$ws.Range("J31").Value = -0.0215191607156528
$ws.Range("K31").Value = 0.0820000861429213
$ws.Range("L31").Value = -0.262428513527976
$ws.Range("M31").Value = 0.804212772961713
$ws.Range("B32").Value = -0.095990521844521
$ws.Range("C32").Value = 0.149892450325444
$ws.Range("D32").Value = -0.640395974821332
$ws.Range("E32").Value = 0.591538664342115
$ws.Range("F32").Value = -0.107996643534211
$ws.Range("G32").Value = 0.147714121272
$ws.Range("H32").Value = -0.731119290452582
$ws.Range("I32").Value = 0.545554411648185
$ws.Range("J32").Value = 0.0537478448708685
$ws.Range("K32").Value = 0.0596219694940898
$ws.Range("L32").Value = 0.901477179082393
$ws.Range("M32").Value = 0.469530470900526
$ws.Range("B33").Value = 0.0799670637759366
$ws.Range("C33").Value = 0.0528337209882084
$ws.Range("D33").Value = 1.5135610795572
$ws.Range("E33").Value = 0.158428548805544
$ws.Range("F33").Value = 0.0795188283101256
$ws.Range("G33").Value = 0.0565996641558021
$ws.Range("H33").Value = 1.40493463161255
$ws.Range("I33").Value = 0.187520779358444
$ws.Range("J33").Value = 0.112362795880062
$ws.Range("K33").Value = 0.0668701858033043
$ws.Range("L33").Value = 1.68031230256443
$ws.Range("M33").Value = 0.123272132992994
$ws.Range("B34").Value = 0.0194863410254025
$ws.Range("C34").Value = 0.0269016618507546
$ws.Range("D34").Value = 0.724354544842215
$ws.Range("E34").Value = 0.488398004520675
$ws.Range("F34").Value = 0.0489120038118365
$ws.Range("G34").Value = 0.0319679956262656
$ws.Range("H34").Value = 1.53003035860182
$ws.Range("I34").Value = 0.16299875638935
$ws.Range("J34").Value = 0.0406623346494248
$ws.Range("K34").Value = 0.0467867715626284
$ws.Range("L34").Value = 0.869098963047589
$ws.Range("M34").Value = 0.411109606997288
$ws.Range("B35").Value = 0.0334024678116172
$ws.Range("C35").Value = 0.0563408134975711
$ws.Range("D35").Value = 0.592864492683572
$ws.Range("E35").Value = 0.563793372982055
$ws.Range("F35").Value = 0.0388406923271813
$ws.Range("G35").Value = 0.0584364950934175
$ws.Range("H35").Value = 0.664664988293531
$ws.Range("I35").Value = 0.518270404339051
$ws.Range("J35").Value = 0.0337430255960825
$ws.Range("K35").Value = 0.0525597301543995
$ws.Range("L35").Value = 0.641993889560676
$ws.Range("M35").Value = 0.534662970093019
$ws.Range("B36").Value = 0.0687141287729896
$ws.Range("C36").Value = 0.0593616291391591
$ws.Range("D36").Value = 1.15755126281837
$ws.Range("E36").Value = 0.273825648477961
$ws.Range("F36").Value = 0.067973256053014
$ws.Range("G36").Value = 0.0759249274439743
$ws.Range("H36").Value = 0.895269292198824
$ws.Range("I36").Value = 0.391607972095222
$ws.Range("J36").Value = 0.14343571311934
$ws.Range("K36").Value = 0.080750457987468
$ws.Range("L36").Value = 1.77628358642375
$ws.Range("M36").Value = 0.109098900419642
$ws.Range("B37").Value = 0.0547680972436089
$ws.Range("C37").Value = 0.101411794805836
$ws.Range("D37").Value = 0.540056483059673
$ws.Range("E37").Value = 0.599776291956376
$ws.Range("F37").Value = 0.14597810046852
$ws.Range("G37").Value = 0.0980659908861985
$ws.Range("H37").Value = 1.48857008581009
$ws.Range("I37").Value = 0.163509865329654
$ws.Range("J37").Value = 0.243239040777967
$ws.Range("K37").Value = 0.10860449915907
$ws.Range("L37").Value = 2.23967738594053
$ws.Range("M37").Value = 0.0494715627225887
$ws.Range("B38").Value = 0.267230567699981
$ws.Range("C38").Value = 0.138398800661861
$ws.Range("D38").Value = 1.93087343547785
$ws.Range("E38").Value = 0.132210633476738
$ws.Range("F38").Value = 0.508295262245764
$ws.Range("G38").Value = 0.135347907926851
$ws.Range("H38").Value = 3.75547188007127
$ws.Range("I38").Value = 0.0204884472467183
$ws.Range("J38").Value = 0.595898656171933
$ws.Range("K38").Value = 0.0914604975601375
$ws.Range("L38").Value = 6.51536643762642
$ws.Range("M38").Value = 0.00401175519270642
$ws.Range("B39").Value = 0.273805684025576
$ws.Range("C39").Value = 0.0236861172500012
$ws.Range("D39").Value = 11.5597538058105
$ws.Range("E39").Value = 0.0063205202742511
$ws.Range("F39").Value = 0.226328465941219
$ws.Range("G39").Value = 0.0445717359839943
$ws.Range("H39").Value = 5.07784722637893
$ws.Range("I39").Value = 0.0371526899737136
$ws.Range("J39").Value = 0.221072875231213
$ws.Range("K39").Value = 0.0477654632651488
$ws.Range("L39").Value = 4.62829961480798
$ws.Range("M39").Value = 0.043460638146096
$ws.Range("B40").Value = 0.316447217012437
$ws.Range("C40").Value = 0.0403591138075913
$ws.Range("D40").Value = 7.84078705298319
$ws.Range("E40").Value = 0.00313470417611607
$ws.Range("F40").Value = 0.272213485377529
$ws.Range("G40").Value = 0.0620898309769207
$ws.Range("H40").Value = 4.38418789509529
$ws.Range("I40").Value = 0.0198781678330264
$ws.Range("J40").Value = 0.27278098654784
$ws.Range("K40").Value = 0.0679045168139014
$ws.Range("L40").Value = 4.01712580174042
$ws.Range("M40").Value = 0.025358980684582
$ws.Range("B41").Value = 0.387580683117061
$ws.Range("C41").Value = 0.117397797587239
$ws.Range("D41").Value = 3.30143061524681
$ws.Range("E41").Value = 0.01149947612315
$ws.Range("F41").Value = 0.328933185899016
$ws.Range("G41").Value = 0.141680734521187
$ws.Range("H41").Value = 2.32165076649731
$ws.Range("I41").Value = 0.0499793474097329
$ws.Range("J41").Value = 0.255396052383797
$ws.Range("K41").Value = 0.154191533725547
$ws.Range("L41").Value = 1.65635587254997
$ws.Range("M41").Value = 0.137707677124421
$ws.Range("B42").Value = 0.284537369067252
$ws.Range("C42").Value = 0.0717808214595598
$ws.Range("D42").Value = 3.9639748233803
$ws.Range("E42").Value = 0.00648482434735788
$ws.Range("F42").Value = 0.224567183708554
$ws.Range("G42").Value = 0.0860772013805728
$ws.Range("H42").Value = 2.60890433363041
$ws.Range("I42").Value = 0.0383479863235434
$ws.Range("J42").Value = 0.172015890286456
$ws.Range("K42").Value = 0.0935321279307803
$ws.Range("L42").Value = 1.83911019766127
$ws.Range("M42").Value = 0.113252273074767
$ws.Range("B43").Value = 0.0579161350800583
$ws.Range("C43").Value = 0.196548204652236
$ws.Range("D43").Value = 0.294666314467397
$ws.Range("E43").Value = 0.775733744457842
$ws.Range("F43").Value = -0.0460868509612041
$ws.Range("G43").Value = 0.2494148669753
$ws.Range("H43").Value = -0.184779887101791
$ws.Range("I43").Value = 0.85799790214702
$ws.Range("J43").Value = -0.143951823498731
$ws.Range("K43").Value = 0.269478252848107
$ws.Range("L43").Value = -0.534187163443835
$ws.Range("M43").Value = 0.607734619339915
$ws.Range("B44").Value = 0.000120363536468705
$ws.Range("C44").Value = 0.000133624684442289
$ws.Range("D44").Value = 0.9007582466598
$ws.Range("E44").Value = 0.397622553794315
$ws.Range("F44").Value = 0.00024200921113486
$ws.Range("G44").Value = 0.000150268837058677
$ws.Range("H44").Value = 1.61050831211504
$ws.Range("I44").Value = 0.151310972274106
$ws.Range("J44").Value = 0.0000238665418758153
$ws.Range("K44").Value = 0.000204407495164094
$ws.Range("L44").Value = 0.116759622031744
$ws.Range("M44").Value = 0.910303337658463
$ws.Range("B45").Value = 0.0090315665997582
$ws.Range("C45").Value = 0.00384670315288063
$ws.Range("D45").Value = 2.34787199344843
$ws.Range("E45").Value = 0.0384588915791509
$ws.Range("F45").Value = 0.00408707565601144
$ws.Range("G45").Value = 0.00426104561490703
$ws.Range("H45").Value = 0.95917200269179
$ws.Range("I45").Value = 0.357971371588739
$ws.Range("J45").Value = -0.00536446199281287
$ws.Range("K45").Value = 0.0050993236676171
$ws.Range("L45").Value = -1.05199480214984
$ws.Range("M45").Value = 0.315187610079838
$ws.Range("B46").Value = 0.00537236281549012
$ws.Range("C46").Value = 0.00296412145942915
$ws.Range("D46").Value = 1.81246379037543
$ws.Range("E46").Value = 0.159047725147206
$ws.Range("F46").Value = 0.0109028961366476
$ws.Range("G46").Value = 0.00554948123827727
$ws.Range("H46").Value = 1.96466942917932
$ws.Range("I46").Value = 0.137176250829159
$ws.Range("J46").Value = 0.0136799390774787
$ws.Range("K46").Value = 0.00488080629456925
$ws.Range("L46").Value = 2.80280311322742
$ws.Range("M46").Value = 0.0613386269628791
$ws.Range("B47").Value = -0.00648817070291021
$ws.Range("C47").Value = 0.011631217496999
$ws.Range("D47").Value = -0.557823865350661
$ws.Range("E47").Value = 0.593959727257267
$ws.Range("F47").Value = -0.011111336231861
$ws.Range("G47").Value = 0.0140623523105523
$ws.Range("H47").Value = -0.790147763793628
$ws.Range("I47").Value = 0.454010135580933
$ws.Range("B48").Value = -0.066892272910515
$ws.Range("C48").Value = 0.0652585520553947
$ws.Range("D48").Value = -1.02503458632875
$ws.Range("E48").Value = 0.326383585935995
$ws.Range("F48").Value = -0.0326580525940371
$ws.Range("G48").Value = 0.0730661511136254
$ws.Range("H48").Value = -0.446965552397175
$ws.Range("I48").Value = 0.663449713405216
$ws.Range("B49").Value = -0.0256513892691586
$ws.Range("C49").Value = 0.0105992679960052
$ws.Range("D49").Value = -2.42010950933841
$ws.Range("E49").Value = 0.0415959285674311
$ws.Range("F49").Value = -0.020645881435336
$ws.Range("G49").Value = 0.0113351080453872
$ws.Range("H49").Value = -1.82141020206135
$ws.Range("I49").Value = 0.108495029256178
$ws.Range("B50").Value = -0.00306185250562964
$ws.Range("C50").Value = 0.00663368141505904
$ws.Range("D50").Value = -0.461561584594485
$ws.Range("E50").Value = 0.652540934014362
$ws.Range("F50").Value = -0.00139552093842605
$ws.Range("G50").Value = 0.00851304304371952
$ws.Range("H50").Value = -0.163927391328721
$ws.Range("I50").Value = 0.872421747845206
$ws.Range("B51").Value = -0.00218397909738343
$ws.Range("C51").Value = 0.000980505380299105
$ws.Range("D51").Value = -2.22740144140484
$ws.Range("E51").Value = 0.242829524751135
$ws.Range("F51").Value = -0.00233257447681332
$ws.Range("G51").Value = 0.00205381318178912
$ws.Range("H51").Value = -1.13572865219482
$ws.Range("I51").Value = 0.441415618368621
$ws.Range("B52").Value = 0.0379341823890264
$ws.Range("C52").Value = 0.016195465674469
$ws.Range("D52").Value = 2.34227179085236
$ws.Range("E52").Value = 0.0435091218745745
$ws.Range("F52").Value = 0.0445921254761612
$ws.Range("G52").Value = 0.0162603070660509
$ws.Range("H52").Value = 2.74239135183756
$ws.Range("I52").Value = 0.0217270733810906
$ws.Range("B53").Value = 0.0254856817725251
$ws.Range("C53").Value = 0.0236826630675467
$ws.Range("D53").Value = 1.07613243070832
$ws.Range("E53").Value = 0.303790734289199
$ws.Range("F53").Value = 0.0468359793881494
$ws.Range("G53").Value = 0.0346815027850382
$ws.Range("H53").Value = 1.35045991745071
$ws.Range("I53").Value = 0.204442893754881
$ws.Range("B54").Value = 0.00000203356807982343
$ws.Range("C54").Value = 0.00000539357948669589
$ws.Range("D54").Value = 0.377034969974865
$ws.Range("E54").Value = 0.713369228606828
$ws.Range("F54").Value = 0.0000021161653878107
$ws.Range("G54").Value = 0.00000736338671086278
$ws.Range("H54").Value = 0.28739022828841
$ws.Range("I54").Value = 0.779124980481246
$ws.Range("B55").Value = 0.0000106647126812373
$ws.Range("C55").Value = 0.00000323554385194765
$ws.Range("D55").Value = 3.29611130902077
$ws.Range("E55").Value = 0.00707885052711643
$ws.Range("F55").Value = 0.0000117054201418346
$ws.Range("G55").Value = 0.00000393719191528211
$ws.Range("H55").Value = 2.97303773697196
$ws.Range("I55").Value = 0.0131188463077965
$ws.Range("B56").Value = -0.393398482943731
$ws.Range("C56").Value = 0.15812401329859
$ws.Range("D56").Value = -2.48791106889543
$ws.Range("E56").Value = 0.0553229216797683
$ws.Range("B57").Value = -0.0947247772606969
$ws.Range("C57").Value = 0.0442000922442614
$ws.Range("D57").Value = -2.14309003558686
$ws.Range("E57").Value = 0.0606500400845647
$ws.Range("B58").Value = -0.159350287644514
$ws.Range("C58").Value = 0.0991120521697175
$ws.Range("D58").Value = -1.60777911622338
$ws.Range("E58").Value = 0.132660379615
$ws.Range("B59").Value = 0.283521878091263
$ws.Range("C59").Value = 0.100141782129387
$ws.Range("D59").Value = 2.83120463868858
$ws.Range("E59").Value = 0.0210939303546739
$ws.Range("B60").Value = -0.553041455344118
$ws.Range("C60").Value = 0.101540501090368
$ws.Range("D60").Value = -5.44651099221902
$ws.Range("E60").Value = 0.00381842464565197
